$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous category grid before writing the new layout
$ws.Range("A1:J19").ClearContents()

$ws.Range("A1").Value = "hus"
$ws.Range("B1").Value = "alabu"
$ws.Range("C1").Value = "leje"
$ws.Range("D1").Value = "boliglån"
$ws.Range("E1").Value = "el"
$ws.Range("F1").Value = "vand"
$ws.Range("G1").Value = "varme"
$ws.Range("H1").Value = "olie"
$ws.Range("I1").Value = "seas"
$ws.Range("J1").Value = "energi nord"

$ws.Range("A2").Value = "forsikring"
$ws.Range("B2").Value = "nykredit a/s"
$ws.Range("C2").Value = "ulykke"
$ws.Range("D2").Value = "indbo"
$ws.Range("E2").Value = "husforsikring"
$ws.Range("F2").Value = "lønsikring"

$ws.Range("A3").Value = "bil"
$ws.Range("B3").Value = "vægtafgift"
$ws.Range("C3").Value = "billån"
$ws.Range("D3").Value = "audi"
$ws.Range("E3").Value = "bilforsikring"

$ws.Range("A4").Value = "lån"
$ws.Range("B4").Value = "gæld"
$ws.Range("C4").Value = "kommune"

$ws.Range("A5").Value = "opsparing"

$ws.Range("A6").Value = "akasse"
$ws.Range("B6").Value = "a-kasse"
$ws.Range("C6").Value = "ftfa"
$ws.Range("D6").Value = "ase"
$ws.Range("E6").Value = "ida"

$ws.Range("A7").Value = "skat"
$ws.Range("B7").Value = "afgift"

$ws.Range("A8").Value = "avis"

$ws.Range("A9").Value = "spotify"

$ws.Range("A10").Value = "netflix"

$ws.Range("A11").Value = "tv"

$ws.Range("A12").Value = "internet"

$ws.Range("A13").Value = "mobil"

$ws.Range("A14").Value = "fitness"

$ws.Range("A15").Value = "transport"
$ws.Range("B15").Value = "rejse"

$ws.Range("A16").Value = "dr-licens"

$ws.Range("B16").Select()
